$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.105
$ws.Range("E2").Value = 0.11635
$ws.Range("F2").Value = 0.006999999999999999
$ws.Range("I2").Value = 0.001177587722684887
$ws.Range("J2").Value = 0.000959752250139885
$ws.Range("K2").Value = 3025.9
$ws.Range("L2").Value = 0.2635962123125975
$ws.Range("M2").Value = 534.273
$ws.Range("N2").Value = 0.01157083733088463
$ws.Range("O2").Value = 0.1765666413298523
$ws.Range("P2").Value = 296.723
$ws.Range("Q2").Value = 0.006426178312084047
$ws.Range("R2").Value = 0.09806107273868932
$ws.Range("S2").Value = 237.55
$ws.Range("T2").Value = 0.4446228800631887
$ws.Range("U2").Value = 13207.4
$ws.Range("V2").Value = 0.2860348117234554
$ws.Range("W2").Value = 0.1206134650376917
$ws.Range("X2").Value = 0.07357343801793793
$ws.Range("Y2").Value = 0.04704002701975378
$ws.Range("Z2").Value = 0.1235499211682527
$ws.Range("AB2").Value = 0.04244073081189914
$ws.Range("AC2").Value = -0.04244073081189914
$ws.Range("AD2").Value = 82296.59999999999
$ws.Range("AE2").Value = 494.9105862749169
$ws.Range("AF2").Value = 82791.51058627492
$ws.Range("AG2").Value = 69584.11058627492
$ws.Range("AH2").Value = 0.6419657939035568
$ws.Range("AI2").Value = 0.7374991116936478
$ws.Range("AJ2").Value = 0.6011159833402374
$ws.Range("AK2").Value = 0.7024979016100468
$ws.Range("AN2").Value = 731.5253333333333
$ws.Range("AP2").Value = 618.5254274335549
$ws.Range("B3").Value = "Grupo Financiero Inbursa, S.A.B. de C.V. (BMV:GFINBUR O)"
$ws.Range("D3").Value = 0.266
$ws.Range("E3").Value = -0.0139
$ws.Range("F3").Value = 0.0517
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 591.6
$ws.Range("L3").Value = 0.4681490860172509
$ws.Range("M3").Value = 34.3
$ws.Range("N3").Value = 0.00519594625301077
$ws.Range("O3").Value = 0.05797836375929682
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("S3").Value = 34.3
$ws.Range("T3").Value = 1
$ws.Range("U3").Value = 2780.2
$ws.Range("V3").Value = 0.4211594685895202
$ws.Range("W3").Value = 0.08106330501507263
$ws.Range("X3").Value = 0.0601118431136857
$ws.Range("Y3").Value = 0.02095146190138693
$ws.Range("Z3").Value = 0.1021452358627825
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.04182423572113846
$ws.Range("AC3").Value = -0.04182423572113846
$ws.Range("AD3").Value = 6523.3
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 6523.3
$ws.Range("AG3").Value = 3743.1
$ws.Range("AH3").Value = 0.4970284808679883
$ws.Range("AI3").Value = 0.4839745967682104
$ws.Range("AJ3").Value = 0.3618479563821971
$ws.Range("AK3").Value = 0.3498747476258132
$ws.Range("B4").Value = "Banco del Bajío, S.A., Institución de Banca Múltiple (BMV:BBAJIO O)"
$ws.Range("D4").Value = 0.148
$ws.Range("E4").Value = 0.229
$ws.Range("F4").Value = -0.145
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 185.6
$ws.Range("L4").Value = 0.3443413729128015
$ws.Range("M4").Value = 9.9
$ws.Range("N4").Value = 0.006068779501011464
$ws.Range("O4").Value = 0.05334051724137932
$ws.Range("P4").Value = 0.02
$ws.Range("Q4").Value = [double]"1.226016060810397e-05"
$ws.Range("R4").Value = 0.0001077586206896552
$ws.Range("S4").Value = 9.880000000000001
$ws.Range("T4").Value = 0.9979797979797981
$ws.Range("U4").Value = 1013.3
$ws.Range("V4").Value = 0.6211610372095875
$ws.Range("W4").Value = 0.1206134650376917
$ws.Range("X4").Value = 0.07357343801793793
$ws.Range("Y4").Value = 0.04704002701975378
$ws.Range("Z4").Value = 0.1755054328537053
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.04243773322277765
$ws.Range("AC4").Value = -0.04243773322277765
$ws.Range("AD4").Value = 2656.5
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 2656.5
$ws.Range("AG4").Value = 1643.2
$ws.Range("AH4").Value = 0.6195484864032837
$ws.Range("AI4").Value = 0.6325602438327459
$ws.Range("AJ4").Value = 0.5018170713085968
$ws.Range("AK4").Value = 0.5157078743370053
$ws.Range("D5").Value = 0.105
$ws.Range("E5").Value = 0.16
$ws.Range("F5").Value = 0.0509
$ws.Range("K5").Value = 1519.8
$ws.Range("L5").Value = 0.3932212160413971
$ws.Range("M5").Value = 97.203
$ws.Range("N5").Value = 0.006101998154390855
$ws.Range("O5").Value = 0.06395775759968417
$ws.Range("P5").Value = 0.003
$ws.Range("Q5").Value = [double]"1.883274637940451e-07"
$ws.Range("R5").Value = [double]"1.973943939992104e-06"
$ws.Range("S5").Value = 97.2
$ws.Range("T5").Value = 0.9999691367550384
$ws.Range("U5").Value = 4533.3
$ws.Range("V5").Value = 0.2845816305391815
$ws.Range("W5").Value = 0.1618151231873256
$ws.Range("X5").Value = 0.07366053112839516
$ws.Range("Y5").Value = 0.08815459205893048
$ws.Range("Z5").Value = 0.1157546062246927
$ws.Range("AB5").Value = 0.04244073081189914
$ws.Range("AC5").Value = -0.04244073081189914
$ws.Range("AD5").Value = 26006.8
$ws.Range("AF5").Value = 26006.8
$ws.Range("AG5").Value = 21473.5
$ws.Range("AH5").Value = 0.6201471272042254
$ws.Range("AI5").Value = 0.7253730510696454
$ws.Range("AJ5").Value = 0.5741086324164778
$ws.Range("AK5").Value = 0.6856227869360179
$ws.Range("B6").Value = "Grupo Elektra, S.A.B. de C.V. (BMV:ELEKTRA *)"
$ws.Range("D6").Value = 0.0839
$ws.Range("K6").Value = -158.5
$ws.Range("L6").Value = -0.05628151409701016
$ws.Range("M6").Value = 143.8
$ws.Range("N6").Value = 0.009526143899094421
$ws.Range("O6").Value = -0.9072555205047319
$ws.Range("P6").Value = 49.8
$ws.Range("Q6").Value = 0.003299040098573728
$ws.Range("R6").Value = -0.3141955835962145
$ws.Range("S6").Value = 94.00000000000001
$ws.Range("T6").Value = 0.6536856745479833
$ws.Range("U6").Value = 1098.2
$ws.Range("V6").Value = 0.07275112120991302
$ws.Range("W6").Value = -0.03218666233449761
$ws.Range("X6").Value = 0.05419445864936823
$ws.Range("Y6").Value = -0.08638112098386583
$ws.Range("Z6").Value = 0.2014045827731213
$ws.Range("AB6").Value = 0.04453361481025302
$ws.Range("AC6").Value = -0.04453361481025302
$ws.Range("AD6").Value = 10668.4
$ws.Range("AF6").Value = 10668.4
$ws.Range("AG6").Value = 9570.199999999999
$ws.Range("AH6").Value = 0.414086486024911
$ws.Range("AI6").Value = 0.7145229994374046
$ws.Range("AJ6").Value = 0.3879994324055867
$ws.Range("AK6").Value = 0.6918583635759004
$ws.Range("B7").Value = "Banco Santander México, S.A., Institución de Banca Múltiple, Grupo Financiero Santander México (BMV:BSMX B)"
$ws.Range("D7").Value = 0.0959
$ws.Range("E7").Value = 0.0727
$ws.Range("F7").Value = -0.0369
$ws.Range("I7").Value = 0.004512880665359091
$ws.Range("J7").Value = 0.003350300380164928
$ws.Range("K7").Value = 887.4
$ws.Range("L7").Value = 0.2962542565266742
$ws.Range("M7").Value = 249.07
$ws.Range("N7").Value = 0.03601098821658353
$ws.Range("O7").Value = 0.2806738787469011
$ws.Range("P7").Value = 246.9
$ws.Range("Q7").Value = 0.03569724571676426
$ws.Range("R7").Value = 0.2782285327924273
$ws.Range("S7").Value = 2.169999999999987
$ws.Range("T7").Value = 0.008712410165816789
$ws.Range("U7").Value = 3782.4
$ws.Range("V7").Value = 0.546866189546736
$ws.Range("W7").Value = 0.1254222435797775
$ws.Range("X7").Value = 0.1516152713683387
$ws.Range("Y7").Value = -0.02619302778856122
$ws.Range("Z7").Value = 0.09952450390257668
$ws.Range("AA7").Value = 0.0003334369832605285
$ws.Range("AB7").Value = 0.04744435523022269
$ws.Range("AC7").Value = -0.04711091824696216
$ws.Range("AD7").Value = 36441.6
$ws.Range("AE7").Value = 494.9105862749169
$ws.Range("AF7").Value = 36936.51058627492
$ws.Range("AG7").Value = 33154.11058627492
$ws.Range("AH7").Value = 0.8422799277054716
$ws.Range("AI7").Value = 0.843341484239622
$ws.Range("AJ7").Value = 0.827392198451574
$ws.Range("AK7").Value = 0.8285335599591875
$ws.Range("AN7").Value = 323.9253333333333
$ws.Range("AP7").Value = 294.7032052113326

$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("F6").ClearContents()

Write-Output "applied"
